$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column L (Num F) values for rows 2-11 -- these were placeholder 0s,
# now filled in with the actual pre-test respondent counts.
$ws.Range("L2").Value = 45
$ws.Range("L3").Value = 100
$ws.Range("L4").Value = 80
$ws.Range("L5").Value = 122
$ws.Range("L6").Value = 88
$ws.Range("L7").Value = 163
$ws.Range("L8").Value = 52
$ws.Range("L9").Value = 155
$ws.Range("L10").Value = 43
$ws.Range("L11").Value = 167

# Add new row 12 (2020 111) data
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "2020 111"
$ws.Range("C12").Value = 12.33333333333333
$ws.Range("D12").Value = 0.7072544488881022
$ws.Range("E12").Value = 14.89473684210526
$ws.Range("F12").Value = 0.8012694386798697
$ws.Range("G12").Value = 2.56140350877193
$ws.Range("H12").Value = 0.1449851042701092
$ws.Range("I12").Value = 0.02304780461542826
$ws.Range("J12").Value = 0.6914341384628476
$ws.Range("K12").Value = 1.068757020484248
$ws.Range("L12").Value = 57
$ws.Range("M12").Value = 14.76576576576577
$ws.Range("N12").Value = 0.5200256647977225
$ws.Range("O12").Value = 17.96396396396396
$ws.Range("P12").Value = 0.6985333796916776
$ws.Range("Q12").Value = 3.198198198198197
$ws.Range("R12").Value = 0.2099349497338852
$ws.Range("S12").Value = 0.01909700426514009
$ws.Range("T12").Value = 0.572910127954203
$ws.Range("U12").Value = 0.8708476184682317
$ws.Range("V12").Value = 111

# Add new row 13 (2020 112) data
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "2020 112"
$ws.Range("C13").Value = 9.376744186046512
$ws.Range("D13").Value = 0.2736919344088321
$ws.Range("E13").Value = 13.64651162790698
$ws.Range("F13").Value = 0.3863891347290506
$ws.Range("G13").Value = 4.269767441860465
$ws.Range("H13").Value = 0.2070365358592693
$ws.Range("I13").Value = 0.0103739092627643
$ws.Range("J13").Value = 0.3112172778829289
$ws.Range("K13").Value = 0.473501677290707
$ws.Range("L13").Value = 215
$ws.Range("M13").Value = 12.24691358024691
$ws.Range("N13").Value = 0.5413548573186054
$ws.Range("O13").Value = 15.75308641975309
$ws.Range("P13").Value = 0.8201796822477849
$ws.Range("Q13").Value = 3.506172839506174
$ws.Range("R13").Value = 0.1974965229485397
$ws.Range("S13").Value = 0.02474924627020214
$ws.Range("T13").Value = 0.7424773881060645
$ws.Range("U13").Value = 0.9827307834369111
$ws.Range("V13").Value = 81

# Match the bordered/bold/centered style used by the rest of column A
# (row index cells) by copying the format from the row above.
$ws.Range("A11").Copy()
$ws.Range("A12").PasteSpecial(-4122)
$ws.Range("A11").Copy()
$ws.Range("A13").PasteSpecial(-4122)

$excel.CutCopyMode = 0
